$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all date values in column A (rows 2-97) forward by one day
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2() + 1
}

# New wind production values for rows 2-30 (column B)
$newValues = @{
    2  = 268
    3  = 0
    4  = 248
    5  = 257
    6  = 288
    7  = 320
    8  = 333
    9  = 364
    10 = 396
    11 = 391
    12 = 388
    13 = 385
    14 = 403
    15 = 433
    16 = 466
    17 = 460
    18 = 450
    19 = 470
    20 = 524
    21 = 611
    22 = 733
    23 = 790
    24 = 829
    25 = 856
    26 = 930
    27 = 1044
    28 = 1120
    29 = 1215
    30 = 1290
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $newValues[$r]
}
